$wb = $excel.ActiveWorkbook

# --- Update Rest Assured sheet C4/C5 payload text (extra _x000d_ token added per line) ---
$wsRest = $wb.Worksheets.Item("Rest Assured")
$c4Value = @"
{_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "first_name": "Ali",_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "last_name": "Ahmad",_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "email": "ali.ahmad@gmail.com",_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "password": "12345",_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "confirm_password": "12345"_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
}
"@
$wsRest.Range("C4").Value = $c4Value

$c5Value = @"
{_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "email": "ali.ahmad@gmail.com",_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
    "password": "12345"_x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d__x000d_
}
"@
$wsRest.Range("C5").Value = $c5Value

# --- Update "After Life - Search Results" sheet rows 3-12 (reshuffled search results) ---
$wsSearch = $wb.Worksheets.Item("After Life - Search Results")
$wsSearch.Range("A3").Value = 'After Life | Netflix Official Site'
$wsSearch.Range("B3").Value = 'https://www.netflix.com/title/80998491'
$wsSearch.Range("A4").Value = 'After Life (TV Series 2019– ) - IMDb'
$wsSearch.Range("B4").Value = 'https://www.imdb.com/title/tt8398600/'
$wsSearch.Range("A5").Value = 'After.Life (2009) - IMDb'
$wsSearch.Range("B5").Value = 'https://www.imdb.com/title/tt0838247/'
$wsSearch.Range("A6").Value = ''
$wsSearch.Range("B6").Value = 'https://www.hitc.com/en-gb/2020/04/26/after-life-dog-brandy-belong-to-ricky-gervais-netflix-anti/'
$wsSearch.Range("A7").Value = ''
$wsSearch.Range("B7").Value = 'https://www.denofgeek.com/tv/after-life-season-3-netflix-renewed/'
$wsSearch.Range("A8").Value = ''
$wsSearch.Range("B8").Value = 'https://en.wikipedia.org/wiki/After_Life_(TV_series)'
$wsSearch.Range("A9").Value = ''
$wsSearch.Range("B9").Value = 'https://www.rogerebert.com/reviews/afterlife-2010'
$wsSearch.Range("A10").Value = 'After Life: Season 1 - Rotten Tomatoes'
$wsSearch.Range("B10").Value = 'https://www.rottentomatoes.com/tv/after_life/s01'
$wsSearch.Range("A11").Value = 'After Life - Rotten Tomatoes'
$wsSearch.Range("B11").Value = 'https://www.rottentomatoes.com/tv/after_life'
$wsSearch.Range("A12").Value = 'Ricky Gervais calls After Life the best thing he''s done. This is ...'
$wsSearch.Range("B12").Value = 'https://www.theguardian.com/tv-and-radio/2019/mar/11/ricky-gervais-calls-after-life-the-best-thing-hes-done-this-is-patently-false'